$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.069.28'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.67%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.777.26'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.12%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.48%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.83%  '

# Row 6: 'USDC' -> 'USDC'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.49%  '

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3795'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.78%  '

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3403'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.36%  '

# Row 9: 'OKB' -> 'OKB'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.97'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.58%  '

# Row 10: 'Polygon' -> 'Polygon'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.183'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.00%  '

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07410'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.75%  '

# Row 12: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.40%  '

# Row 13: 'Solana' -> 'Solana'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.50'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.75%  '

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.392'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.34%  '

# Row 15: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.778.80'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.94%  '

# Row 16: 'Chainlink' -> 'Chainlink'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.041'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.41%  '

# Row 17: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001080'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.13%  '

# Row 18: 'TRON' -> 'TRON'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06654'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.01%  '

# Row 19: 'Litecoin' -> 'Litecoin'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.11'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.77%  '

# Row 20: 'Dai' -> 'Dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.39%  '

# Row 21: 'Uniswap' -> 'Uniswap'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.534'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.04%  '

# Row 22: 'Avalanche' -> 'Avalanche'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.21'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.00%  '

# Row 23: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.069.38'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.63%  '

# Row 24: 'Cosmos' -> 'Cosmos'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.15'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -7.54%  '

# Row 25: 'Toncoin' -> 'Toncoin'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.379'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.49%  '

# Row 26: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.496'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -6.60%  '

# Row 27: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.452'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.97%  '

# Row 28: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.92'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -5.37%  '

# Row 29: 'Monero' -> 'Monero'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.29'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.59%  '

# Row 30: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.979.85'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.90%  '

# Row 31: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.42'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.19%  '

# Row 32: 'HuobiToken' -> 'HuobiToken'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.003'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.89%  '

# Row 33: 'Filecoin' -> 'Filecoin'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.970'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -6.16%  '

# Row 34: 'Stellar' -> 'Stellar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08641'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.90%  '

# Row 35: 'Aptos' -> 'Aptos'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.97'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.48%  '

# Row 36: 'WEMIXTOKEN' -> 'WEMIXTOKEN'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.628'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.03%  '

# Row 37: 'TheSandbox' -> 'InternetComputer(DFINITY)'
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.353'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.12%  '

# Row 38: 'InternetComputer(DFINITY)' -> 'TheSandbox'
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6777'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.43%  '

# Row 39: 'Hedera' -> 'Hedera'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06264'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.29%  '

# Row 40: 'FraxShare' -> 'Algorand'
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2166'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.82%  '

# Row 41: 'Algorand' -> 'VeChain'
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02312'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.57%  '

# Row 42: 'VeChain' -> 'FraxShare'
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.550'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.18%  '

# Row 43: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.226'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.36%  '

# Row 44: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.19'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.44%  '

# Row 45: 'Frax' -> 'Frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.005'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.49%  '

# Row 46: 'Decentraland' -> 'Decentraland'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6364'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.19%  '

# Row 47: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.849'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.07%  '

# Row 48: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.117'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.95%  '

# Row 49: 'Quant' -> 'Quant'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.41'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.82%  '

# Row 50: 'Cronos' -> 'Cronos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07097'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.56%  '

# Row 51: 'Aave' -> 'Aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.29'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.98%  '
